# Insertar descuento para cada cliente en el archivo Excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the client discount amount on the invoice (row 18, "Client discount").
# This flows through the existing formulas (AMOUNT, SUBTOTAL, TAX, TOTAL)
# automatically on recalculation.
$ws.Range("E18").Value = 120

# Update the footer placeholder with the RPA developer's contact details
$ws.Range("A31").Value = "RPA Developer - RPADeveloper@Uipath.com"
